$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '27.379.23'
$ws.Range("E2").Value = '  +0.24%  '
Set-TextValue $ws.Range("D3") '1.715.99'
$ws.Range("E3").Value = '  +0.25%  '
Set-TextValue $ws.Range("D4") '1.008'
$ws.Range("E4").Value = '  +0.40%  '
Set-TextValue $ws.Range("D5") '224.76'
$ws.Range("E5").Value = '  +0.32%  '
Set-TextValue $ws.Range("D6") '0.5279'
$ws.Range("E6").Value = '  -0.33%  '
Set-TextValue $ws.Range("D7") '1.008'
$ws.Range("E7").Value = '  +0.30%  '
Set-TextValue $ws.Range("D8") '0.06663'
$ws.Range("E8").Value = '  +1.84%  '
Set-TextValue $ws.Range("D9") '0.2648'
$ws.Range("E9").Value = '  +0.41%  '
Set-TextValue $ws.Range("D10") '20.79'
$ws.Range("E10").Value = '  -1.06%  '
Set-TextValue $ws.Range("D11") '0.07751'
$ws.Range("E11").Value = '  +1.42%  '
Set-TextValue $ws.Range("D12") '4.470'
$ws.Range("E12").Value = '  -2.10%  '
Set-TextValue $ws.Range("D13") '1.952.15'
$ws.Range("E13").Value = '  +0.27%  '
Set-TextValue $ws.Range("D14") '1.717.05'
$ws.Range("E14").Value = '  +0.28%  '
Set-TextValue $ws.Range("D15") '0.5793'
$ws.Range("E15").Value = '  +0.82%  '
Set-TextValue $ws.Range("D16") '0.0₅8206'
$ws.Range("E16").Value = '  +0.23%  '
Set-TextValue $ws.Range("D17") '67.80'
$ws.Range("E17").Value = '  +0.84%  '
Set-TextValue $ws.Range("D18") '27.384.01'
$ws.Range("E18").Value = '  +0.26%  '
Set-TextValue $ws.Range("D19") '219.78'
$ws.Range("E19").Value = '  +1.58%  '
Set-TextValue $ws.Range("D20") '1.007'
$ws.Range("E20").Value = '  +0.23%  '
Set-TextValue $ws.Range("D21") '4.648'
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("E22").Value = '  -0.36%  '
Set-TextValue $ws.Range("D23") '6.055'
$ws.Range("E23").Value = '  +1.75%  '
$ws.Range("E24").Value = '  +0.28%  '
Set-TextValue $ws.Range("D25") '145.41'
$ws.Range("E25").Value = '  +2.01%  '
Set-TextValue $ws.Range("D26") '1.724'
$ws.Range("E26").Value = '  -0.46%  '
Set-TextValue $ws.Range("D27") '0.1205'
$ws.Range("E27").Value = '  -1.21%  '
Set-TextValue $ws.Range("D28") '7.225'
$ws.Range("E28").Value = '  -0.61%  '
Set-TextValue $ws.Range("D29") '16.18'
$ws.Range("E29").Value = '  -0.64%  '
Set-TextValue $ws.Range("D30") '0.05333'
$ws.Range("E30").Value = '  -0.87%  '
Set-TextValue $ws.Range("D31") '1.295'
$ws.Range("E31").Value = '  +0.38%  '
Set-TextValue $ws.Range("D32") '3.481'
$ws.Range("E32").Value = '  -0.38%  '
Set-TextValue $ws.Range("D33") '3.366'
$ws.Range("E33").Value = '  -1.36%  '
Set-TextValue $ws.Range("D34") '1.641'
$ws.Range("E34").Value = '  +0.03%  '
Set-TextValue $ws.Range("D35") '2.836'
$ws.Range("E35").Value = '  -1.10%  '
Set-TextValue $ws.Range("D36") '0.9544'
$ws.Range("E36").Value = '  +0.33%  '
Set-TextValue $ws.Range("D37") '2.402'
$ws.Range("E37").Value = '  -0.91%  '
Set-TextValue $ws.Range("D38") '0.5885'
$ws.Range("E38").Value = '  +0.39%  '
Set-TextValue $ws.Range("D39") '1.187.52'
$ws.Range("E39").Value = '  +14.29%  '
Set-TextValue $ws.Range("D40") '0.01652'
$ws.Range("E40").Value = '  +1.33%  '
Set-TextValue $ws.Range("D41") '5.814'
$ws.Range("E41").Value = '  -1.09%  '
$ws.Range("E42").Value = '  +0.28%  '
Set-TextValue $ws.Range("D43") '0.8428'
$ws.Range("E43").Value = '  +0.40%  '
Set-TextValue $ws.Range("D44") '101.29'
$ws.Range("E44").Value = '  +0.17%  '
Set-TextValue $ws.Range("D45") '1.859.18'
$ws.Range("E45").Value = '  +0.23%  '
$ws.Range("E46").Value = '  +2.65%  '
$ws.Range("E47").Value = '  -0.82%  '
Set-TextValue $ws.Range("D48") '0.4558'
$ws.Range("E48").Value = '  +1.49%  '
Set-TextValue $ws.Range("D49") '1.006'
$ws.Range("E49").Value = '  +0.35%  '
Set-TextValue $ws.Range("D50") '8.134'
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("E51").Value = '  -0.09%  '
